# Refresh Price (column D) and Volume(1h) (column E) figures for each
# cryptocurrency row, matching the latest scrape from coinranking.com.
# Figures are stored as plain text in the sheet, so numeric-looking values
# are written with a leading apostrophe (Excel's 'store as text' marker) and
# the cell style is reset to Normal afterwards so no stray number formatting
# is introduced.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = '64.655.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +2.93%  '
$ws.Range("E2").Style = "Normal"

# Row 3: Ethereum
$ws.Range("D3").Value = '2.525.43'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.54%  '
$ws.Range("E3").Style = "Normal"

# Row 4: TetherUSD
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("E4").Style = "Normal"

# Row 5: BNB
$ws.Range("D5").Value = '''579.81'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("E5").Style = "Normal"

# Row 6: Solana
$ws.Range("D6").Value = '''153.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +4.10%  '
$ws.Range("E6").Style = "Normal"

# Row 7: USDC
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E7").Style = "Normal"

# Row 8: XRP
$ws.Range("E8").Value = '  +0.87%  '
$ws.Range("E8").Style = "Normal"

# Row 9: LidoStakedEther
$ws.Range("D9").Value = '2.525.81'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.57%  '
$ws.Range("E9").Style = "Normal"

# Row 10: Dogecoin
$ws.Range("D10").Value = '''0.113'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.64%  '
$ws.Range("E10").Style = "Normal"

# Row 11: TRON
$ws.Range("E11").Value = '  -1.68%  '
$ws.Range("E11").Style = "Normal"

# Row 12: Toncoin
$ws.Range("E12").Value = '  +0.02%  '
$ws.Range("E12").Style = "Normal"

# Row 13: Cardano
$ws.Range("D13").Value = '''0.355'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("E13").Style = "Normal"

# Row 14: Avalanche
$ws.Range("E14").Value = '  +0.77%  '
$ws.Range("E14").Style = "Normal"

# Row 15: ShibaInu
$ws.Range("E15").Value = '  +1.31%  '
$ws.Range("E15").Style = "Normal"

# Row 16: WrappedliquidstakedEther2.0
$ws.Range("D16").Value = '2.982.11'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.48%  '
$ws.Range("E16").Style = "Normal"

# Row 17: WrappedBTC
$ws.Range("D17").Value = '64.251.52'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +2.35%  '
$ws.Range("E17").Style = "Normal"

# Row 18: WrappedEther
$ws.Range("D18").Value = '2.527.97'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.95%  '
$ws.Range("E18").Style = "Normal"

# Row 19: Uniswap
$ws.Range("D19").Value = '''7.88'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.70%  '
$ws.Range("E19").Style = "Normal"

# Row 20: Chainlink
$ws.Range("D20").Value = '''10.97'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.41%  '
$ws.Range("E20").Style = "Normal"

# Row 21: Polkadot
$ws.Range("D21").Value = '''4.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.23%  '
$ws.Range("E21").Style = "Normal"

# Row 22: BitcoinCash
$ws.Range("D22").Value = '''329.40'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.82%  '
$ws.Range("E22").Style = "Normal"

# Row 23: SuiNetwork
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E23").Style = "Normal"

# Row 24: Dai
$ws.Range("E24").Value = '  +0.05%  '
$ws.Range("E24").Style = "Normal"

# Row 25: Aptos
$ws.Range("D25").Value = '''10.07'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.07%  '
$ws.Range("E25").Style = "Normal"

# Row 26: Litecoin
$ws.Range("D26").Value = '''65.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E26").Style = "Normal"

# Row 27: Bittensor
$ws.Range("D27").Value = '''645.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.06%  '
$ws.Range("E27").Style = "Normal"

# Row 28: PEPE
$ws.Range("D28").Value = '''0.0000105'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.31%  '
$ws.Range("E28").Style = "Normal"

# Row 30: Fetch.AI
$ws.Range("E30").Value = '  +4.83%  '
$ws.Range("E30").Style = "Normal"

# Row 31: Binance-PegBSC-USD
$ws.Range("D31").Value = '''0.995'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.18%  '
$ws.Range("E31").Style = "Normal"

# Row 32: InternetComputer(DFINITY)
$ws.Range("D32").Value = '''8.05'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.85%  '
$ws.Range("E32").Style = "Normal"

# Row 33: PancakeSwap
$ws.Range("E33").Value = '  +1.60%  '
$ws.Range("E33").Style = "Normal"

# Row 34: Kaspa
$ws.Range("D34").Value = '''0.137'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.97%  '
$ws.Range("E34").Style = "Normal"

# Row 35: FirstDigitalUSD
$ws.Range("D35").Value = '''0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.07%  '
$ws.Range("E35").Style = "Normal"

# Row 36: ImmutableX
$ws.Range("E36").Value = '  +0.78%  '
$ws.Range("E36").Style = "Normal"

# Row 37: NEARProtocol
$ws.Range("D37").Value = '''4.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("E37").Style = "Normal"

# Row 38: RenderToken
$ws.Range("D38").Value = '''5.55'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.22%  '
$ws.Range("E38").Style = "Normal"

# Row 39: Monero
$ws.Range("D39").Value = '''154.37'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.45%  '
$ws.Range("E39").Style = "Normal"

# Row 40: EthereumClassic
$ws.Range("D40").Value = '''18.96'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.38%  '
$ws.Range("E40").Style = "Normal"

# Row 41: PolygonEcosystemToken
$ws.Range("E41").Value = '  +0.89%  '
$ws.Range("E41").Style = "Normal"

# Row 42: dogwifhat
$ws.Range("D42").Value = '''2.83'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.18%  '
$ws.Range("E42").Style = "Normal"

# Row 43: Stacks
$ws.Range("D43").Value = '''1.80'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.38%  '
$ws.Range("E43").Style = "Normal"

# Row 44: Aave
$ws.Range("D44").Value = '''163.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +7.16%  '
$ws.Range("E44").Style = "Normal"

# Row 45: USDe
$ws.Range("E45").Value = '  +0.02%  '
$ws.Range("E45").Style = "Normal"

# Row 46: BabyDogeCoin
$ws.Range("D46").Value = '0.0₆0301'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("E46").Style = "Normal"

# Row 47: WhiteBITCoin
$ws.Range("E47").Value = '  +1.73%  '
$ws.Range("E47").Style = "Normal"

# Row 48: Filecoin
$ws.Range("E48").Value = '  +1.55%  '
$ws.Range("E48").Style = "Normal"

# Row 49: InjectiveProtocol
$ws.Range("D49").Value = '''21.40'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +4.36%  '
$ws.Range("E49").Style = "Normal"

# Row 50: Mantle
$ws.Range("E50").Value = '  +2.03%  '
$ws.Range("E50").Style = "Normal"

# Row 51: Hedera
$ws.Range("D51").Value = '''0.0518'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
$ws.Range("E51").Style = "Normal"
